$wb = $excel.ActiveWorkbook

$wsUser       = $wb.Worksheets.Item(1)   # User
$wsSubscriber = $wb.Worksheets.Item(2)   # Subscriber
$wsAddress    = $wb.Worksheets.Item(3)   # Address
$wsSponsor    = $wb.Worksheets.Item(4)   # Sponsor
$wsPicture    = $wb.Worksheets.Item(5)   # Picture

# --- User sheet: add "Sample" column (C) with sample data for every mapped field ---
$wsUser.Range("C1").Value = "Sample"
$wsUser.Range("C2").Value = "c7f5a5a8-6648-443d-b8cf-490efc32adb6"   # Id / Guid
$wsUser.Range("C3").Value = "John"                                   # FirstName
$wsUser.Range("C4").Value = "Smith"                                  # LastName
$wsUser.Range("C5").Value = "jsmith@yahoo.com"                       # Email
$wsUser.Hyperlinks.Add($wsUser.Range("C5"), "mailto:jsmith@yahoo.com")
$wsUser.Range("C6").Value = "904-555-9273"                           # Phone
$wsUser.Range("C7").Value = "ashid7hdeaip78ai"                       # hashedPassword
$wsUser.Range("C8").Value = "30f8c3e1-0186-4330-8d4c-0b1ed44cc403"   # CreateUser / id

$wsUser.Range("C9").Value = 42719                                    # CreateDate / DateTime
$wsUser.Range("C9").NumberFormat = "mm-dd-yy"

$wsUser.Range("C10").Value = "30f8c3e1-0186-4330-8d4c-0b1ed44cc403"  # UpdateUser / id

$wsUser.Range("C11").Value = 42719                                   # UpdateDate / DateTime
$wsUser.Range("C9").Copy() | Out-Null
$wsUser.Range("C11").PasteSpecial(-4122) | Out-Null

$wsUser.Range("C12").Value = 3                                       # Version / Int

$wsUser.Columns.Item(3).AutoFit() | Out-Null

# --- Remaining sheets: just add the "Sample" header in column C ---
$wsSubscriber.Range("C1").Value = "Sample"
$wsAddress.Range("C1").Value = "Sample"
$wsSponsor.Range("C1").Value = "Sample"
$wsPicture.Range("C1").Value = "Sample"

# --- Restore selections / active sheet to match the saved workbook state ---
$wsUser.Range("C17").Select() | Out-Null
$wsSubscriber.Range("C1").Select() | Out-Null
$wsAddress.Range("C1").Select() | Out-Null
$wsSponsor.Range("C1").Select() | Out-Null
$wsPicture.Range("C1").Select() | Out-Null
